$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(770).Insert()

$ws.Cells.Item(770, 1).Value = 10
$ws.Cells.Item(770, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(770, 3).Value = 'La Araucanía'
$ws.Cells.Item(770, 4).Value = 45131
$ws.Cells.Item(770, 5).Value = 9
$ws.Cells.Item(770, 6).Value = 100112045
$ws.Cells.Item(770, 7).Value = 'Zapallo'
$ws.Cells.Item(770, 8).Value = 'Camote'
$ws.Cells.Item(770, 9).Value = '1a (guarda)'
$ws.Cells.Item(770, 10).Value = 850
$ws.Cells.Item(770, 11).Value = 500
$ws.Cells.Item(770, 12).Value = 500
$ws.Cells.Item(770, 13).Value = 500
$ws.Cells.Item(770, 14).Value = '$/kilo (volumen en unidades)'
$ws.Cells.Item(770, 15).Value = 'Región del Maule'
$ws.Cells.Item(770, 16).Value = 500
$ws.Cells.Item(770, 17).Value = 1
$ws.Cells.Item(770, 18).Value = 'Hortaliza'
